$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.421699999999998
$ws.Range("B7").Value = 5.4102
$ws.Range("A9").Value = -21.75380000000001
$ws.Range("B12").Value = 5.335699999999997
$ws.Range("B14").Value = 6.3125
$ws.Range("D15").Value = -8.66
$ws.Range("A18").Value = -22.10920000000001
$ws.Range("A20").Value = -19.02829999999999
$ws.Range("B26").Value = 4.078200000000004
$ws.Range("A27").Value = -21.75019999999997
$ws.Range("B27").Value = 5.033100000000002
$ws.Range("B29").Value = 5.124200000000001
$ws.Range("D33").Value = -7.744499999999999
$ws.Range("A35").Value = -20.2456
$ws.Range("D35").Value = -7.6343
$ws.Range("B37").Value = 8.929599999999999
$ws.Range("B38").Value = 4.739100000000002
$ws.Range("D38").Value = -8.567900000000002
$ws.Range("D43").Value = -8.225799999999998
$ws.Range("D44").Value = -7.714700000000001
$ws.Range("D47").Value = -7.502800000000001
$ws.Range("B51").Value = 5.6953
$ws.Range("D51").Value = -7.509699999999992
$ws.Range("B52").Value = 5.108300000000003
$ws.Range("B55").Value = 5.143799999999999
$ws.Range("D57").Value = -8.262800000000002
$ws.Range("D63").Value = -8.029799999999996
$ws.Range("A69").Value = -21.70069999999999
$ws.Range("B69").Value = 5.371499999999998
$ws.Range("B70").Value = 6.977300000000001
$ws.Range("D70").Value = -7.199800000000002
$ws.Range("A76").Value = -19.28859999999999
$ws.Range("A78").Value = -19.93189999999998
$ws.Range("B81").Value = 5.579200000000002
$ws.Range("A82").Value = -21.8768
$ws.Range("A83").Value = -21.9241
$ws.Range("B83").Value = 6.347700000000008
$ws.Range("D88").Value = -8.246099999999998
$ws.Range("A93").Value = -20.53719999999997
$ws.Range("D99").Value = -7.451299999999999
$ws.Range("B102").Value = 8.735000000000007
